$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("D2").Value = '28.129.94'
$ws.Range("E2").Value = '  -0.28%  '

# Row 3
$ws.Range("D3").Value = '1.873.61'
$ws.Range("E3").Value = '  -1.88%  '

# Row 4
$ws.Range("E4").Value = '  +0.31%  '

# Row 5
$ws.Range("D5").Value = "'314.13"
$ws.Range("E5").Value = '  -0.19%  '

# Row 6
$ws.Range("E6").Value = '  +0.28%  '

# Row 7
$ws.Range("D7").Value = "'0.5058"
$ws.Range("E7").Value = '  +0.09%  '

# Row 8
$ws.Range("D8").Value = "'0.3837"
$ws.Range("E8").Value = '  -2.37%  '

# Row 9
$ws.Range("D9").Value = "'0.08578"
$ws.Range("E9").Value = '  -8.00%  '

# Row 10
$ws.Range("D10").Value = "'1.116"
$ws.Range("E10").Value = '  -2.41%  '

# Row 11
$ws.Range("E11").Value = '  -0.65%  '

# Row 12
$ws.Range("D12").Value = "'6.327"
$ws.Range("E12").Value = '  -1.14%  '

# Row 13
$ws.Range("B13").Value = 'WrappedEther'
$ws.Range("C13").Value = 'https://coinranking.com/coin/Mtfb0obXVh59u+wrappedether-weth'
$ws.Range("D13").Value = '1.894.79'
$ws.Range("E13").Value = '  -0.65%  '

# Row 14
$ws.Range("B14").Value = 'Solana'
$ws.Range("C14").Value = 'https://coinranking.com/coin/zNZHO_Sjf+solana-sol'
$ws.Range("D14").Value = "'20.66"
$ws.Range("E14").Value = '  -1.08%  '

# Row 15
$ws.Range("D15").Value = "'1.004"
$ws.Range("E15").Value = '  +0.33%  '

# Row 16
$ws.Range("D16").Value = "'7.160"
$ws.Range("E16").Value = '  -2.09%  '

# Row 17
$ws.Range("D17").Value = "'0.00001101"
$ws.Range("E17").Value = '  -2.06%  '

# Row 18
$ws.Range("D18").Value = "'91.14"

# Row 19
$ws.Range("D19").Value = "'0.06612"
$ws.Range("E19").Value = '  -0.04%  '

# Row 20
$ws.Range("D20").Value = "'18.16"
$ws.Range("E20").Value = '  +0.88%  '

# Row 21
$ws.Range("D21").Value = "'1.003"
$ws.Range("E21").Value = '  +0.28%  '

# Row 22
$ws.Range("E22").Value = '  -1.98%  '

# Row 23
$ws.Range("D23").Value = '28.158.63'
$ws.Range("E23").Value = '  -0.36%  '

# Row 24
$ws.Range("D24").Value = "'11.39"
$ws.Range("E24").Value = '  -0.59%  '

# Row 25
$ws.Range("D25").Value = "'2.276"
$ws.Range("E25").Value = '  -1.91%  '

# Row 26
$ws.Range("D26").Value = '2.107.95'
$ws.Range("E26").Value = '  -0.88%  '

# Row 27
$ws.Range("D27").Value = "'2.556"
$ws.Range("E27").Value = '  -2.10%  '

# Row 28
$ws.Range("D28").Value = "'157.47"
$ws.Range("E28").Value = '  -0.49%  '

# Row 29
$ws.Range("D29").Value = "'20.73"

# Row 30
$ws.Range("D30").Value = "'127.13"
$ws.Range("E30").Value = '  -0.19%  '

# Row 31
$ws.Range("D31").Value = "'0.1050"
$ws.Range("E31").Value = '  -2.27%  '

# Row 32
$ws.Range("D32").Value = "'1.058"
$ws.Range("E32").Value = '  -4.19%  '

# Row 33
$ws.Range("D33").Value = "'5.593"
$ws.Range("E33").Value = '  -0.96%  '

# Row 34
$ws.Range("D34").Value = "'3.606"
$ws.Range("E34").Value = '  -0.15%  '

# Row 35
$ws.Range("D35").Value = "'9.663"
$ws.Range("E35").Value = '  -0.19%  '

# Row 36
$ws.Range("D36").Value = "'0.02438"
$ws.Range("E36").Value = '  +0.68%  '

# Row 37
$ws.Range("D37").Value = "'0.06550"
$ws.Range("E37").Value = '  -1.67%  '

# Row 38
$ws.Range("D38").Value = "'0.2169"
$ws.Range("E38").Value = '  -1.03%  '

# Row 39
$ws.Range("D39").Value = "'1.203"
$ws.Range("E39").Value = '  -3.34%  '

# Row 40
$ws.Range("D40").Value = "'1.243"
$ws.Range("E40").Value = '  -5.15%  '

# Row 41
$ws.Range("D41").Value = "'11.53"
$ws.Range("E41").Value = '  +0.26%  '

# Row 42
$ws.Range("D42").Value = "'0.6367"
$ws.Range("E42").Value = '  -1.38%  '

# Row 43
$ws.Range("E43").Value = '  -2.30%  '

# Row 44
$ws.Range("B44").Value = 'Frax'
$ws.Range("C44").Value = 'https://coinranking.com/coin/KfWtaeV1W+frax-frax'
$ws.Range("D44").Value = "'1.004"
$ws.Range("E44").Value = '  +0.37%  '

# Row 45
$ws.Range("B45").Value = 'EnergySwap'
$ws.Range("C45").Value = 'https://coinranking.com/coin/SbWqqTui-+energyswap-ens'
$ws.Range("D45").Value = "'13.16"
$ws.Range("E45").Value = '  -1.16%  '

# Row 46
$ws.Range("B46").Value = 'Decentraland'
$ws.Range("C46").Value = 'https://coinranking.com/coin/tEf7-dnwV3BXS+decentraland-mana'
$ws.Range("D46").Value = "'0.5981"
$ws.Range("E46").Value = '  -0.91%  '

# Row 47
$ws.Range("D47").Value = "'3.677"
$ws.Range("E47").Value = '  -1.14%  '

# Row 48
$ws.Range("B48").Value = 'WEMIXTOKEN'
$ws.Range("C48").Value = 'https://coinranking.com/coin/08CsQa-Ov+wemixtoken-wemix'
$ws.Range("D48").Value = "'1.278"
$ws.Range("E48").Value = '  -0.42%  '

# Row 49
$ws.Range("B49").Value = 'EOS'
$ws.Range("C49").Value = 'https://coinranking.com/coin/iAzbfXiBBKkR6+eos-eos'
$ws.Range("D49").Value = "'1.231"
$ws.Range("E49").Value = '  +3.80%  '

# Row 50
$ws.Range("B50").Value = 'NEARProtocol'
$ws.Range("C50").Value = 'https://coinranking.com/coin/DCrsaMv68+nearprotocol-near'
$ws.Range("D50").Value = "'1.985"
$ws.Range("E50").Value = '  -1.88%  '

# Row 51
$ws.Range("B51").Value = 'Quant'
$ws.Range("C51").Value = 'https://coinranking.com/coin/bauj_21eYVwso+quant-qnt'
$ws.Range("D51").Value = "'121.13"
$ws.Range("E51").Value = '  -1.61%  '
